$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '87.320.77'
$ws.Range("E2").Value = '  -2.07%  '

$ws.Range("D3").Value = '3.169.19'
$ws.Range("E3").Value = '  -6.36%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.31'
$ws.Range("E5").Value = '  -7.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '611.66'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.378'
$ws.Range("E7").Value = '  -11.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.664'
$ws.Range("E8").Value = '  -0.42%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").Value = '3.161.91'
$ws.Range("E10").Value = '  -6.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.534'
$ws.Range("E11").Value = '  -16.27%  '

$ws.Range("E12").Value = '  +4.23%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000241'
$ws.Range("E13").Value = '  -17.53%  '

$ws.Range("D14").Value = '3.759.96'
$ws.Range("E14").Value = '  -6.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.25'
$ws.Range("E15").Value = '  -7.32%  '

$ws.Range("D16").Value = '87.080.52'
$ws.Range("E16").Value = '  -2.17%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '32.00'
$ws.Range("E17").Value = '  -14.16%  '

$ws.Range("D18").Value = '3.158.32'
$ws.Range("E18").Value = '  -6.36%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.02'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.40'
$ws.Range("E20").Value = '  -10.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '414.49'
$ws.Range("E21").Value = '  -11.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.47'
$ws.Range("E22").Value = '  -13.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.07'
$ws.Range("E23").Value = '  -11.66%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.17'
$ws.Range("E24").Value = '  -8.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.89'
$ws.Range("E25").Value = '  -7.92%  '

$ws.Range("D26").Value = '3.338.82'
$ws.Range("E26").Value = '  -6.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '73.43'
$ws.Range("E27").Value = '  -9.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000129'
$ws.Range("E28").Value = '  -11.41%  '

$ws.Range("E29").Value = '  -0.12%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.159'
$ws.Range("E30").Value = '  -18.67%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  -0.18%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '543.81'
$ws.Range("E32").Value = '  -8.61%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.22'
$ws.Range("E33").Value = '  -13.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.29'
$ws.Range("E34").Value = '  -17.92%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.72'
$ws.Range("E35").Value = '  -8.86%  '

$ws.Range("E36").Value = '  -13.53%  '

$ws.Range("E37").Value = '  -9.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '21.77'
$ws.Range("E38").Value = '  -9.55%  '

$ws.Range("B39").Value = 'WhiteBITCoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '21.77'
$ws.Range("E39").Value = '  -0.36%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  +0.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.97'
$ws.Range("E41").Value = '  -8.35%  '

$ws.Range("E42").Value = '  +0.01%  '

$ws.Range("E43").Value = '  -12.11%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.369'
$ws.Range("E44").Value = '  -16.23%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.59'
$ws.Range("E45").Value = '  -6.36%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '172.75'
$ws.Range("E46").Value = '  -9.08%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.19'
$ws.Range("E47").Value = '  -8.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.126'
$ws.Range("E48").Value = '  -1.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.23'
$ws.Range("E49").Value = '  -15.49%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.95'
$ws.Range("E50").Value = '  -13.64%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.696'
$ws.Range("E51").Value = '  -12.50%  '
